$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '30.333.43'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'" + '  -2.88%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'" + '1.941.36'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'" + '  -2.92%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'" + '  +0.19%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'" + '250.76'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'" + '  -2.62%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'" + '0.7184'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'" + '  -8.42%  '
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'" + '0.9999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'" + '  +0.11%  '
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'" + '0.3364'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'" + '  -4.81%  '
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'" + '28.79'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'" + '  -1.48%  '
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'" + '0.07313'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'" + '  +3.75%  '
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'" + '0.8183'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'" + '  -5.91%  '
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'" + '0.08142'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'" + '  -0.80%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'" + '1.938.68'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'" + '  -3.07%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'" + '5.528'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'" + '  -1.23%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'" + '95.43'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'" + '  -5.44%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'" + '14.88'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'" + '  -3.90%  '
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'" + '30.352.41'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'" + '  -2.83%  '
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'" + '0.000008256'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'" + '  +3.73%  '
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'" + '254.32'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'" + '  -7.45%  '
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'" + '5.897'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'" + '  -1.10%  '
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'" + '2.193.75'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'" + '  -2.77%  '
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'" + '  +0.06%  '
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'" + '1.001'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'" + '  +0.32%  '
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'" + '6.973'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'" + '  -2.47%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'" + '9.870'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'" + '  -2.42%  '
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'" + '160.35'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'" + '  -2.67%  '
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'" + '2.455'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'" + '  +3.65%  '
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'" + '  -2.81%  '
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'" + '0.1320'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'" + '  -11.57%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'" + '1.570'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'" + '  -2.61%  '
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'" + '1.345'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'" + '  -1.25%  '
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'" + '4.493'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'" + '  -2.84%  '
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'" + '4.247'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'" + '  -4.57%  '
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'" + '0.05261'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'" + '  +0.64%  '
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'" + '1.273'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'" + '  +3.13%  '
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'" + '0.7546'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'" + '  -2.97%  '
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'" + '2.737'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'" + '  -2.62%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'" + '0.01990'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'" + '  -1.12%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'" + '2.834'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'" + '  -2.58%  '
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'" + '82.16'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'" + '  +2.88%  '
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'" + '6.565'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'" + '  -2.69%  '
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'" + '0.4590'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'" + '  -3.01%  '
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'" + '2.028'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'" + '  -5.75%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'" + '0.8465'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'" + '  -0.72%  '
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'" + '1.000'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'" + '  +0.10%  '
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'" + '102.67'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'" + '  -3.63%  '
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'" + '9.914'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'" + '  -0.49%  '
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'" + '7.463'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'" + '  -3.72%  '
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'" + '37.11'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'" + '  +0.68%  '
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'" + '0.4212'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'" + '  -3.07%  '
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'" + '1.510'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'" + '  -0.06%  '
$ws.Range("E51").ClearFormats()
